# Update the cryptos price/volume table with the latest scraped values.
# Each value is prefixed with a literal apostrophe so Excel always stores it
# as text (preserving exact formatting such as trailing zeros, thousands
# separators written as dots, leading/trailing spaces and the %% sign)
# instead of re-interpreting look-alike numbers and dropping formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''55.735.98'
$ws.Range('E2').Value = '''  +6.50%  '
$ws.Range('D3').Value = '''2.505.27'
$ws.Range('E3').Value = '''  +8.18%  '
$ws.Range('E4').Value = '''  -0.13%  '
$ws.Range('D5').Value = '''490.29'
$ws.Range('E5').Value = '''  +12.56%  '
$ws.Range('D6').Value = '''140.79'
$ws.Range('E6').Value = '''  +16.77%  '
$ws.Range('D7').Value = '''1.00'
$ws.Range('E7').Value = '''  +0.38%  '
$ws.Range('D8').Value = '''0.514'
$ws.Range('E8').Value = '''  +10.43%  '
$ws.Range('D9').Value = '''2.495.94'
$ws.Range('E9').Value = '''  +7.87%  '
$ws.Range('D10').Value = '''0.0987'
$ws.Range('E10').Value = '''  +12.22%  '
$ws.Range('D11').Value = '''5.53'
$ws.Range('E11').Value = '''  +6.27%  '
$ws.Range('D12').Value = '''0.332'
$ws.Range('E12').Value = '''  +9.62%  '
$ws.Range('E13').Value = '''  +1.96%  '
$ws.Range('D14').Value = '''2.937.12'
$ws.Range('E14').Value = '''  +8.68%  '
$ws.Range('D15').Value = '''55.689.21'
$ws.Range('E15').Value = '''  +6.25%  '
$ws.Range('E16').Value = '''  +10.86%  '
$ws.Range('E17').Value = '''  +16.89%  '
$ws.Range('D18').Value = '''2.502.64'
$ws.Range('E18').Value = '''  +7.66%  '
$ws.Range('E19').Value = '''  +11.97%  '
$ws.Range('D20').Value = '''323.16'
$ws.Range('E20').Value = '''  +8.32%  '
$ws.Range('D21').Value = '''10.03'
$ws.Range('E21').Value = '''  +12.45%  '
$ws.Range('D22').Value = '''0.999'
$ws.Range('E22').Value = '''  +0.09%  '
$ws.Range('D23').Value = '''5.77'
$ws.Range('E23').Value = '''  +12.03%  '
$ws.Range('D24').Value = '''58.33'
$ws.Range('E24').Value = '''  +9.13%  '
$ws.Range('E25').Value = '''  +16.52%  '
$ws.Range('D26').Value = '''0.413'
$ws.Range('E26').Value = '''  +13.60%  '
$ws.Range('E27').Value = '''  +0.75%  '
$ws.Range('D28').Value = '''2.622.29'
$ws.Range('E28').Value = '''  +8.70%  '
$ws.Range('D29').Value = '''7.48'
$ws.Range('E29').Value = '''  +8.05%  '
$ws.Range('D30').Value = '''0.0₃0797'
$ws.Range('E30').Value = '''  +18.82%  '
$ws.Range('E31').Value = '''  +0.55%  '
$ws.Range('D32').Value = '''150.20'
$ws.Range('E32').Value = '''  +4.68%  '
$ws.Range('D33').Value = '''18.27'
$ws.Range('E33').Value = '''  +7.77%  '
$ws.Range('E34').Value = '''  +12.78%  '
$ws.Range('D35').Value = '''5.23'
$ws.Range('E35').Value = '''  +10.84%  '
$ws.Range('D36').Value = '''0.877'
$ws.Range('E36').Value = '''  +7.27%  '
$ws.Range('D37').Value = '''3.71'
$ws.Range('E37').Value = '''  +6.57%  '
$ws.Range('D38').Value = '''1.12'
$ws.Range('E38').Value = '''  +12.81%  '
$ws.Range('E39').Value = '''  +7.52%  '
$ws.Range('D40').Value = '''0.0558'
$ws.Range('E40').Value = '''  +11.46%  '
$ws.Range('B41').Value = '''FirstDigitalUSD'
$ws.Range('C41').Value = '''https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D41').Value = '''0.997'
$ws.Range('E41').Value = '''  +0.18%  '
$ws.Range('B42').Value = '''Mantle'
$ws.Range('C42').Value = '''https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D42').Value = '''0.611'
$ws.Range('E42').Value = '''  +17.30%  '
$ws.Range('D43').Value = '''3.44'
$ws.Range('E43').Value = '''  +9.29%  '
$ws.Range('E44').Value = '''  +9.46%  '
$ws.Range('D45').Value = '''4.75'
$ws.Range('E45').Value = '''  +20.15%  '
$ws.Range('D46').Value = '''2.003.89'
$ws.Range('E46').Value = '''  +5.30%  '
$ws.Range('D47').Value = '''0.0917'
$ws.Range('E47').Value = '''  +11.44%  '
$ws.Range('B48').Value = '''Bittensor'
$ws.Range('C48').Value = '''https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '''256.23'
$ws.Range('E48').Value = '''  +35.84%  '
$ws.Range('B49').Value = '''WhiteBITCoin'
$ws.Range('C49').Value = '''https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D49').Value = '''10.12'
$ws.Range('E49').Value = '''  -0.38%  '
$ws.Range('E50').Value = '''  +9.71%  '
$ws.Range('E51').Value = '''  +12.12%  '
